# Vega Modelo de Temuco - Pepino ensalada: add a new weekly price record.
# Insert a new row at 329 (pushes the existing rows 329:347 down to 330:348)
# and populate it with the new observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(329).Insert()

$ws.Cells.Item(329, 1).Value = 10
$ws.Cells.Item(329, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(329, 3).Value = "La Araucanía"
$ws.Cells.Item(329, 4).Value = 44516
$ws.Cells.Item(329, 5).Value = 9
$ws.Cells.Item(329, 6).Value = 100112043
$ws.Cells.Item(329, 7).Value = "Pepino ensalada"
$ws.Cells.Item(329, 8).Value = "Sin especificar"
$ws.Cells.Item(329, 9).Value = "Primera"
$ws.Cells.Item(329, 10).Value = 235
$ws.Cells.Item(329, 11).Value = 8000
$ws.Cells.Item(329, 12).Value = 9000
$ws.Cells.Item(329, 13).Value = 8468
$ws.Cells.Item(329, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(329, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(329, 16).Value = 141
$ws.Cells.Item(329, 17).Value = 60
$ws.Cells.Item(329, 18).Value = "Hortaliza"
